$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet view / selection changes ---
# Move the active selection from F14 to D9 (this also clears the scrolled
# "topLeftCell" (was A5) that pinned the view).
$ws.Range("D9").Select()

# --- Data changes ---
# reag1 number of wells: 48 -> 96
$ws.Range("D7").Value = 96

# Clear the stray "#" (s49) marker text from column A on several rows,
# leaving the cell formatting untouched.
$ws.Range("A10").ClearContents()
$ws.Range("A11").ClearContents()
$ws.Range("A12").ClearContents()
$ws.Range("A13").ClearContents()
$ws.Range("A14").ClearContents()
$ws.Range("A15").ClearContents()
$ws.Range("A16").ClearContents()
$ws.Range("A31").ClearContents()
$ws.Range("A53").ClearContents()
$ws.Range("A57").ClearContents()

# chem3_abbreviation: AcNH3I -> PyrrolidiniumIodide
# (leading apostrophe preserves the cell's existing quote-prefix / forced-text style)
$ws.Range("D27").Value = "'PyrrolidiniumIodide"

# chem1_abbreviation (lower bound): GBL -> DMSO
$ws.Range("D21").Value = "DMSO"

# Reagent concentration values
$ws.Range("D44").Value = 2.2799999999999998
$ws.Range("D45").Value = 2.85
$ws.Range("D51").Value = 4.18
